# This script replaces the "Test setup" section (heading + specs) and the
# old "Foer test" stub with the new, expanded content:
#   - "Test setup" heading becomes a single run "Test setup" tagged en-US
#   - CPU/RAM/OS lines are preserved as-is (still en-US tagged)
#   - "HDD: SSD" loses its en-US language tag
#   - The old empty paragraph + "Foer test" fragment (en-US tagged, with
#     spell-check wrapping around "Foer") is replaced by a new "Foer test"
#     heading (no language tag) followed by a new paragraph of explanatory
#     text about using jvisualvm, and two trailing empty paragraphs.

$d = $word.ActiveDocument

# Locate the "Test setup" heading paragraph (style Overskrift1, text starting
# with "Test") so we don't depend on hard-coded character offsets.
$startIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t.StartsWith("Test ")) {
        $startIndex = $i
        break
    }
}

$startPos = $d.Paragraphs.Item($startIndex).Range.Start
$endPos = $d.Content.End

# Remove the old "Test setup" ... "Foer test" block entirely (including all
# its paragraph marks), then splice in the replacement content as raw OOXML
# so the exact paragraph/run/proofErr structure can be controlled precisely.
$rng = $d.Range($startPos, $endPos)
$rng.Delete()

$insertXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Overskrift1"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Test setup</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">CPU: </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Intel® Core™ i7-2640M CPU </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>@2,80GHz</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>RAM: 4GB</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>OS: Windows 7 Professional SP1 – 64bit</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>HDD: SSD</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Overskrift1"/></w:pPr><w:r><w:t>Før test</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>T</w:t></w:r><w:r><w:t>il</w:t></w:r><w:r><w:t xml:space="preserve"> denne test anvendes </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>jvisualvm</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> til at beregne hvor meget RAM og CPU</w:t></w:r><w:r><w:t xml:space="preserve"> systemet anvender. Samtidig vil der blive evalueret på hvordan ressourceforbruget er og om der allokeres nogle ressourcer. Da systemet er i hvile bør det være muligt at allokere meget få ressourcer løbende. Hvis systemet er designet til lange perioder med ingen aktivitet</w:t></w:r><w:r><w:t>, og ressourcerne i hvile skal minimeres, vil det være muligt at allokere de nødvendige ressourcer én gang, i stedet for at oprette og nedlægge dem. Det vil dog være helt op til hvad der er målet for performance.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0"/></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0"/></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target = $d.Range($startPos, $startPos)
$target.InsertXML($insertXml)
